# Update "想去人数" (F column) counts across the 展览, 演出, and 全部类型
# sheets to reflect the latest generated output (gh-pages update).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 838
$ws1.Range("F3").Value = 568
$ws1.Range("F6").Value = 1160
$ws1.Range("F8").Value = 57
$ws1.Range("F11").Value = 1221
$ws1.Range("F14").Value = 905
$ws1.Range("F18").Value = 76
$ws1.Range("F20").Value = 807
$ws1.Range("F21").Value = 1752
$ws1.Range("F22").Value = 3142
$ws1.Range("F23").Value = 919
$ws1.Range("F24").Value = 91
$ws1.Range("F25").Value = 2302
$ws1.Range("F26").Value = 671
$ws1.Range("F27").Value = 9
$ws1.Range("F28").Value = 3163
$ws1.Range("F29").Value = 651
$ws1.Range("F30").Value = 797
$ws1.Range("F33").Value = 744
$ws1.Range("F34").Value = 149
$ws1.Range("F35").Value = 141
$ws1.Range("F36").Value = 49
$ws1.Range("F38").Value = 1127
$ws1.Range("F40").Value = 413
$ws1.Range("F43").Value = 208
$ws1.Range("F44").Value = 138
$ws1.Range("F46").Value = 56

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 94

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 838
$ws4.Range("F3").Value = 568
$ws4.Range("F5").Value = 1160
$ws4.Range("F8").Value = 1221
$ws4.Range("F10").Value = 905
$ws4.Range("F16").Value = 76
$ws4.Range("F17").Value = 807
$ws4.Range("F18").Value = 1752
$ws4.Range("F19").Value = 3142
$ws4.Range("F20").Value = 919
$ws4.Range("F21").Value = 91
$ws4.Range("F23").Value = 2302
$ws4.Range("F24").Value = 9
$ws4.Range("F25").Value = 3163
$ws4.Range("F26").Value = 651
$ws4.Range("F27").Value = 797
$ws4.Range("F34").Value = 94
$ws4.Range("F35").Value = 744
$ws4.Range("F36").Value = 149
$ws4.Range("F37").Value = 141
$ws4.Range("F38").Value = 49
$ws4.Range("F41").Value = 1127
$ws4.Range("F44").Value = 413
$ws4.Range("F46").Value = 208
$ws4.Range("F47").Value = 138
$ws4.Range("F49").Value = 56

$wb.Save()
